$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$cs = $s.ThemeColorScheme
Write-Host "Count:" $cs.Count
try {
    $cs.Item(1).RGB = 255
    Write-Host "set color ok"
} catch {
    Write-Host "set color failed: $($_.Exception.Message)"
}
